# "Inclusion of error rows in comparison file"
#
# The comparison sheet previously dropped the source DataFrame's index
# column ("Unnamed: 0") and kept a row whose ISSUE AGE failed to parse
# into a "before→after" transition (it was left as a raw number, 34,
# instead of e.g. "74→45"). That row is an error row that should not be
# in the comparison output, and the index column should be carried
# through as column A.
#
# This script:
#   1. Deletes worksheet row 2 (POLID 123456 / the malformed error row).
#   2. Restores the "Unnamed: 0" index header in A1, matching the header
#      styling already used by the other header cells.
#   3. Strips the leftover header-ish style from the index values in
#      column A so the remaining data rows render like normal data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Remove the error row entirely; rows below shift up automatically.
$ws.Rows("2").Delete()

# 2. Give A1 the same look (bold / bordered / centered-top) as the other
#    header cells, then set its text.
$ws.Range("B1").Copy()
$ws.Range("A1").PasteSpecial(-4122)
$ws.Range("A1").Value = "Unnamed: 0"

# 3. The old index column (A2:A3, post-delete) carried the same bold
#    header-like style as A1 used to have; clear it so the numeric index
#    values look like ordinary data cells.
$ws.Range("A2:A3").ClearFormats()
